$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ3 = $wb.Worksheets.Item(2)      # currently "2022-Q3"

# Duplicate the existing quarter sheet. Excel places the new copy right
# after the source sheet and makes it the active sheet (tabSelected).
$wsQ3.Copy($null, $wsQ3)

$wsQ4 = $wb.Worksheets.Item(2)      # original sheet -> becomes "2022-Q4"
$wsQ3New = $wb.Worksheets.Item(3)   # the copy -> stays "2022-Q3" (data untouched)

$wsQ4.Name = "2022-Q4"
$wsQ3New.Name = "2022-Q3"

# --- Update "总计" sheet: row 2 becomes the new Q4 entry, row 3 gets the
#     previous Q3 entry that used to live in row 2. ---
$wsTotal.Range("B2").Value = "2022-Q4"

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.01
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# --- Update the "2022-Q4" sheet with the new quarter's fund data. ---
# Header row + A2 pick up the "总计" sheet's header formatting (style index 2)
# instead of the old quarter-sheet header style (style index 1).
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2").PasteSpecial(-4122)

$wsQ4.Range("A2").Value = 0

# B2, D2, E2, F2, G2 are stored as text in the workbook even though some of
# them look numeric - force text storage the same way Excel does
# (NumberFormat "@" before entry), then drop back to the default "Normal"
# style so no extra style record is left behind attached to the cell.
$wsQ4.Range("B2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "165531"
$wsQ4.Range("B2").Style = "Normal"

$wsQ4.Range("C2").Value = "信诚多策略灵活配置混合（LOF）"

$wsQ4.Range("D2:G2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "0.89"
$wsQ4.Range("E2").Value = "72.25"
$wsQ4.Range("F2").Value = "1.02"
$wsQ4.Range("G2").Value = "0.0091"
$wsQ4.Range("D2:G2").Style = "Normal"

$wsQ4.Range("H2").Value = 9

# Match "总计" sheet's page margins (0.75/0.75/1/1/0.5/0.5 in inches).
$wsQ4.PageSetup.LeftMargin = 54
$wsQ4.PageSetup.RightMargin = 54
$wsQ4.PageSetup.TopMargin = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36
